$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (values are stored as text strings,
# matching the inlineStr/shared-string cell type used in the source workbook).
$updates = @{
    "D2" = "302.47"
    "E2" = "-0.37%"
    "D3" = "37.22"
    "E3" = "7.11%"
    "D4" = "4.994"
    "E4" = "-3.54%"
    "D5" = "0.07811"
    "E5" = "0.30%"
    "D6" = "2.197"
    "E6" = "-4.16%"
    "D7" = "8.008"
    "E7" = "0.08%"
    "D8" = "0.9167"
    "E8" = "-1.27%"
    "D9" = "0.09710"
    "E9" = "-3.61%"
    "D10" = "0.1872"
    "E10" = "3.60%"
    "D11" = "0.08626"
    "E11" = "1.55%"
    "D12" = "0.03551"
    "E12" = "2.26%"
    "D13" = "0.09953"
    "E13" = "0.51%"
    "D14" = "0.001476"
    "E14" = "-0.49%"
    "D15" = "0.005708"
    "E15" = "-2.19%"
    "D16" = "3.460"
    "E16" = "-0.23%"
    "D17" = "4.039"
    "E17" = "1.17%"
    "D18" = "2.389"
    "E18" = "13.37%"
    "E19" = "0.61%"
    "D20" = "0.1309"
    "E20" = "-1.34%"
    "D21" = "4.789"
    "E21" = "5.55%"
    "D22" = "0.2298"
    "E22" = "-1.60%"
    "D23" = "0.04618"
    "E23" = "0.11%"
    "E24" = "0.93%"
    "D25" = "0.004780"
    "E25" = "7.81%"
    "D26" = "0.0001406"
    "E26" = "8.16%"
    "E27" = "39.73%"
    "D39" = "0.01776"
    "E39" = "1.24%"
    "D40" = "0.04738"
    "E40" = "0.51%"
    "D41" = "0.008045"
    "E41" = "4.91%"
    "D42" = "0.1392"
    "E42" = "-1.06%"
    "D43" = "0.007821"
    "E43" = "10.95%"
    "D44" = "0.002098"
    "E44" = "-9.19%"
    "D45" = "0.009949"
    "E45" = "8.37%"
    "D46" = "0.00006205"
    "E46" = "3.61%"
    "D47" = "0.00000000753"
    "E47" = "0.40%"
    "D48" = "7.352"
    "E48" = "169.08%"
    "E49" = "-0.28%"
    "D50" = "0.00002109"
    "E50" = "0.40%"
    "D51" = "0.0002008"
    "E51" = "0.40%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage so numeric-looking strings (e.g. "302.47", "-0.37%")
    # are kept as literal text instead of being parsed into numbers/percentages.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
